$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "bis 9:00 Uhr / 12:00-14:00 / 19:00-21:00" row), shifting all following rows up.
$ws.Rows(2).Delete()

# Restore the view selection as captured in the saved workbook.
$ws.Range("I16").Select()
